# Updated cryptos list on Sun Aug 20 19:57:57 UTC 2023 with GitHub Actions
# Refresh per-coin Price (col D) and Volume(1h) (col E) values, plus the
# Aave / BabyDogeCoin row swap in rows 45-46, to match the new scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.398.09'
$ws.Range("D3").Value = '1.693.61'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("D4").Value = '''1.010'
$ws.Range("D5").Value = '''219.19'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").Value = '''0.5491'
$ws.Range("E6").Value = '  +4.06%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").Value = '''0.2754'
$ws.Range("E8").Value = '  +1.61%  '
$ws.Range("D9").Value = '''0.06459'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("D11").Value = '''0.07678'
$ws.Range("E11").Value = '  +2.44%  '
$ws.Range("D12").Value = '1.685.05'
$ws.Range("E12").Value = '  -2.07%  '
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").Value = '''0.5844'
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").Value = '''0.000008364'
$ws.Range("E15").Value = '  -1.95%  '
$ws.Range("D16").Value = '''65.50'
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").Value = '26.429.77'
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").Value = '''4.929'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '''1.011'
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").Value = '''192.27'
$ws.Range("E21").Value = '  +1.47%  '
$ws.Range("D22").Value = '''6.253'
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("D23").Value = '''1.011'
$ws.Range("D24").Value = '''148.83'
$ws.Range("E24").Value = '  +2.77%  '
$ws.Range("E25").Value = '  +6.87%  '
$ws.Range("D26").Value = '''7.920'
$ws.Range("E26").Value = '  +2.61%  '
$ws.Range("D27").Value = '''15.80'
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("D28").Value = '''0.06296'
$ws.Range("E28").Value = '  -5.45%  '
$ws.Range("D29").Value = '''1.385'
$ws.Range("E29").Value = '  +2.09%  '
$ws.Range("D30").Value = '''1.331'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").Value = '''3.606'
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("D32").Value = '''3.606'
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("D33").Value = '''1.683'
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("D34").Value = '''1.043'
$ws.Range("E34").Value = '  +1.08%  '
$ws.Range("E35").Value = '  -1.58%  '
$ws.Range("D36").Value = '''2.414'
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("D37").Value = '''2.713'
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("D39").Value = '''6.190'
$ws.Range("E39").Value = '  -2.81%  '
$ws.Range("D40").Value = '1.115.97'
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("D41").Value = '''0.8820'
$ws.Range("E41").Value = '  -0.88%  '
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").Value = '''101.60'
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("D44").Value = '1.844.36'
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''57.55'
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '''0.00000000108'
$ws.Range("E46").Value = '  -6.35%  '
$ws.Range("D47").Value = '''8.196'
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").Value = '''1.007'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").Value = '''0.05273'
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("D50").Value = '''6.107'
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").Value = '''0.4304'
$ws.Range("E51").Value = '  +0.05%  '
